$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture the current (pre-edit) contents of rows 511 and 512 (A:R) ---
# These two rows hold the "Primera"/"Segunda" quality records for the most
# recent sampling date in the historical block. A new week's sample
# (date 45265) needs to be inserted above the rest of the history, so the
# current contents of 511/512 are shifted down two rows (to 513/514) while
# 511/512 themselves get the new date but keep their other values.

$colCount = 18  # columns A..R

$row511 = @()
$row512 = @()
for ($c = 1; $c -le $colCount; $c++) {
    $row511 += ,$ws.Cells.Item(511, $c).Value2
    $row512 += ,$ws.Cells.Item(512, $c).Value2
}

# --- New sampling date for the top (most recent) pair of rows ---
$ws.Cells.Item(511, 4).Value2 = 45265
$ws.Cells.Item(512, 4).Value2 = 45265

# --- Make room: insert two new rows after the (now updated) row 512 ---
$ws.Rows("513:514").Insert()

# --- Populate the freshly inserted rows with the captured old 511/512 data ---
for ($c = 1; $c -le $colCount; $c++) {
    $ws.Cells.Item(513, $c).Value2 = $row511[$c - 1]
    $ws.Cells.Item(514, $c).Value2 = $row512[$c - 1]
}
